# Refactor status and feasibility sheet data:
#  - rename "Site status" -> "Status data"
#  - rename "Site feasibility" -> "Feasibility data"
#  - add a new comment at A3 ("...1") on both of those sheets,
#    matching the existing comments already present at B3:F3 / A4:A8

$wb = $excel.ActiveWorkbook

$wsStatus = $wb.Worksheets.Item("Site status")
$wsStatus.Name = "Status data"

$wsFeasibility = $wb.Worksheets.Item("Site feasibility")
$wsFeasibility.Name = "Feasibility data"

$wsStatus.Range("A3").AddComment("...1")
$wsFeasibility.Range("A3").AddComment("...1")
